# BangVangThanhTich.xlsx — commit 0712118
# "Nguyen Le Hoang Dung Tao co so du lieu"
#
# Row 9 (STT=4) gets its date / Ho ten / Thanh tich filled in, the "Ho ten"
# column is widened to fit the new name, and a new date-formatted cell
# style is introduced for the date cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date cell: 1-May (serial 40299) — gets a new numeric date format (d-mmm)
# on top of the existing font/border used throughout the table.
$ws.Range("B9").Value = 40299
$ws.Range("B9").NumberFormat = "d-mmm"

# New name / achievement text for row 9.
$ws.Range("C9").Value = "Nguyễn Lê Hoàng Dũng"
$ws.Range("D9").Value = "Tạo cơ sở dữ liệu"

# Widen the "Ho ten" column so the new name fits (target stored width is
# 26.28515625 "characters"; ColumnWidth only round-trips to whole pixels
# for this column's font, so 25.5 is the closest input that lands on the
# nearest achievable width).
$ws.Columns.Item(3).ColumnWidth = 25.5
